$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 obsolete date columns (3-nov, 6-nov, 7-nov, 8-nov, 9-nov) -> CE:CI
$ws.Range("CE1:CI1").EntireColumn.Delete()

# After the delete, the former CJ:CM (10-nov,13-nov,14-nov,15-nov) shift left to CE:CH.
# Append a new column (CI) for 16-nov: copy formatting from the now-last column (CH) first.
$ws.Range("CH1:CH11").Copy()
$ws.Range("CI1:CI11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("CI1").Value = "16-nov"

$ws.Range("CI2").Value = 11
$ws.Range("CI3").Value = 9
$ws.Range("CI4").Value = 8
$ws.Range("CI5").Value = 10
$ws.Range("CI6").Value = 11
$ws.Range("CI7").Value = 8
$ws.Range("CI8").Value = 13
$ws.Range("CI9").Value = 12
$ws.Range("CI10").Value = 20
$ws.Range("CI11").Value = 0

$ws.Range("CM6").Select()
